$d = $word.ActiveDocument

# 1) Partial Gastrectomy bullet list -> single run, joined with spaces
$text1 = "- Removes the tumor - Does not remove lymph nodes - Best suited for: - Small adenocarcinoma - GI Stromal Tumors"
$d.Content.Find.Execute($text1, $true, $false, $false, $false, $false, $true, 1, $false, $text1, 2) | Out-Null

# 2) Subtotal Gastrectomy bullet list -> single run, joined with spaces
$text2 = "- Removes bottom 2/3 of stomach - Removes nearby lymph nodes - Reconstruction with small intestine"
$d.Content.Find.Execute($text2, $true, $false, $false, $false, $false, $true, 1, $false, $text2, 2) | Out-Null

# 3) Proximal Tumors bullet list -> single run, joined with spaces
$text3 = "- Located near the top of the stomach - Challenging area for surgery"
$d.Content.Find.Execute($text3, $true, $false, $false, $false, $false, $true, 1, $false, $text3, 2) | Out-Null

# 4) Total Gastrectomy bullet list -> single run, joined with spaces
$text4 = "- Removes all of the stomach - Reconstruction with small intestine - Needed for those with CDH1 mutations"
$d.Content.Find.Execute($text4, $true, $false, $false, $false, $false, $true, 1, $false, $text4, 2) | Out-Null

# 5) "esophageal" -> "stomach" in the spread sentence
$d.Content.Find.Execute(
    "Some esophageal cancers can spread inside the abdomen", $true, $false, $false, $false, $false,
    $true, 1, $false, "Some stomach cancers can spread inside the abdomen", 2) | Out-Null

# 6) "esophageal cancer" -> "gastric cancer" in the laparoscopy sentence
$d.Content.Find.Execute(
    "Not all patients with esophageal cancer need a laparoscopy.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Not all patients with gastric cancer need a laparoscopy.", 2) | Out-Null

# 7) Remove the following paragraph entirely:
#    "In general, laparoscopy is considered for cancers that invade from the esophagus into the stomach."
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`n", "`a") -eq "In general, laparoscopy is considered for cancers that invade from the esophagus into the stomach.") {
        $p.Range.Delete()
        break
    }
}
